$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny floating point re-computation of the existing A6 timestamp
# (05:00:14.796 on 2025-08-09) picked up by the latest automatic pull.
$ws.Range("A6").Value2 = 45878.20850458333

# Append the new hourly reading (row 7 - 06:00:12) coming from the
# automatic WSL update.
$ws.Range("A7").NumberFormat = $ws.Range("A6").NumberFormat
$ws.Range("A7").Value2 = 45878.25014258089
$ws.Cells.Item(7, 2).Value2 = 2025
$ws.Cells.Item(7, 3).Value2 = 37
$ws.Cells.Item(7, 4).Value2 = 13.11
$ws.Cells.Item(7, 5).Value2 = 92.65000000000001
$ws.Cells.Item(7, 6).Value2 = 0
$ws.Cells.Item(7, 7).Value2 = 0.53
$ws.Cells.Item(7, 8).Value2 = "WNW"
$ws.Cells.Item(7, 9).Value2 = 0
$ws.Cells.Item(7, 10).Value2 = "06:00:12"
